$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1/J1, copying H1's formatting (style) first, then
# overwriting the value/text so the existing cell style (bold, border,
# centered) carries over to the new header cells.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Fill I2:I17 with 1, and J2:J17 with the same value as H2:H17 (column H is
# unstyled data, so plain value assignment matches).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
